$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-124 down to 28-125
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new record
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44998
$ws.Cells.Item(27, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112038
$ws.Cells.Item(27, 7).Value = "Cebollín baby"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 270
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = 1750
$ws.Cells.Item(27, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 875
$ws.Cells.Item(27, 17).Value = 2
$ws.Cells.Item(27, 18).Value = "Hortaliza"
